# PAS-6576 Update "individual VIN retrieval" logic to use ENTRY DATE and VALID
# Made appropriate changes to the VIN Upload files
# Made some remarks in the tests

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ENTRY DATE (column AI) for the first data row moved forward a year,
# from 2000-01-01 to 2001-01-01.
$ws.Range("AI2").Value = 20010101

# MODEL_TEXT (column F) for rows 3-5 used to share the "Gt" placeholder;
# they are now distinct markers used by the VIN-retrieval tests.
$ws.Range("F3").Value = "invalidVIN"
$ws.Range("F4").Value = "SecondValid"
$ws.Range("F5").Value = "ThirdValid"

# Reset the sheet's saved view/selection state (scroll back to the
# top-left and select J13 instead of the previous AA7/top-left=S1 state).
$ws.Range("J13").Select()
